# Adapt column header formatting to respective input file names.
# - Rename header cells from "<name>_old" / "<name>_new" to
#   "<name>_FV2304" / "<name>_FV2310"
# - Turn the data range into an Excel Table (ListObject) with AutoFilter
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row -------------------------------------------------

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J => "<name>_FV2304"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2304"
}

# Column K stays "diff" (unchanged)
$ws.Cells.Item(1, 11).Value = "diff"

# Columns L-U => "<name>_FV2310"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2310"
}

# --- 2) Convert the used range into a Table with an AutoFilter -----------

$tableRange = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3) Freeze the header row ---------------------------------------------

$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
